$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The primary diagnosis (participants) query in B2 had a bug: it matched
# diagnosis -> participant -> study -> sample in a single chained MATCH,
# which silently drops participants whose study/sample do not also satisfy
# the pattern, and it returned samples via a plain (unsorted) collect().
# Replace it with the fixed query that uses OPTIONAL MATCH so participants
# aren't dropped, and sorts the collected sample ids.
$fixedParticipantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in ['Adrenal Cortical Carcinoma']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

$ws.Range("B2").Value = $fixedParticipantsQuery

# Writing the longer replacement text would otherwise cause the wrapped-text
# row to auto-grow on save; restore the original explicit row height so the
# layout stays exactly as it was.
$ws.Rows.Item(2).RowHeight = 157.5

# Move the active selection from C4 to E4 to match the saved workbook state.
$ws.Range("E4").Select()
